$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Insert two new columns before column D (shifts old D:K -> F:M)
$ws.Columns("D:E").Insert()

# Step 2: Copy number formatting/styles from column F (old D, now shifted) into new D:E
# so the new columns inherit the same per-row formatting (date format row, number format rows, etc.)
$ws.Range("F7:G102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Step 3: Populate the new D and E columns with their data values.
# Each entry: "row,Dvalue,Evalue" where value is a number, the literal NA (maps to the
# shared "NA" text already used elsewhere in the sheet) or BLANK (cell left empty).
$rowData = @(
    "7,43465,43373",
    "8,19900,18800",
    "9,NA,NA",
    "10,NA,NA",
    "11,BLANK,BLANK",
    "12,NA,NA",
    "13,0,0",
    "14,0,0",
    "15,-700,-600",
    "16,BLANK,BLANK",
    "17,5700,5200",
    "18,14200,13600",
    "19,BLANK,BLANK",
    "20,-10000,-9100",
    "21,4800,5200",
    "22,0,0",
    "23,4100,4600",
    "24,1100,500",
    "25,0,0",
    "26,3000,4000",
    "27,3000,4000",
    "28,0,0",
    "29,300,0",
    "30,0,0",
    "31,0,0",
    "32,10000,9100",
    "33,3300,4000",
    "34,0,0",
    "35,3300,4000",
    "38,43465,43373",
    "39,BLANK,BLANK",
    "40,BLANK,BLANK",
    "41,15900,21200",
    "42,14800,16300",
    "43,0,0",
    "44,0,0",
    "45,0,0",
    "46,0,0",
    "47,0,0",
    "48,40200,39800",
    "49,19800,19900",
    "50,0,0",
    "51,0,0",
    "52,1100,1800",
    "53,0,0",
    "54,1786500,1735300",
    "55,BLANK,BLANK",
    "56,BLANK,BLANK",
    "57,0,0",
    "58,0,0",
    "59,9900,11200",
    "60,0,0",
    "61,24100,24000",
    "62,0,0",
    "63,0,0",
    "64,0,0",
    "65,0,0",
    "66,1604200,1556900",
    "67,BLANK,BLANK",
    "68,0,0",
    "69,0,0",
    "70,0,0",
    "71,0,0",
    "72,175900,174200",
    "73,0,0",
    "74,0,0",
    "75,0,0",
    "76,182300,178400",
    "77,0,0",
    "80,43465,43373",
    "81,3300,4000",
    "82,BLANK,BLANK",
    "83,700,600",
    "84,0,0",
    "85,0,0",
    "86,0,0",
    "87,0,0",
    "88,0,0",
    "89,4300,3900",
    "90,BLANK,BLANK",
    "91,-1100,-1100",
    "92,0,0",
    "93,0,0",
    "94,-58500,-49400",
    "95,BLANK,BLANK",
    "96,-400,-400",
    "97,0,0",
    "98,0,0",
    "99,0,0",
    "100,46600,35400",
    "101,0,0",
    "102,-7600,-10000"
)

foreach ($entry in $rowData) {
    $parts = $entry.Split(",")
    $r = [int]$parts[0]
    $dVal = $parts[1]
    $eVal = $parts[2]

    if ($dVal -ne "BLANK") {
        if ($dVal -eq "NA") {
            $ws.Cells.Item($r, 4).Value = "NA"
        } else {
            $ws.Cells.Item($r, 4).Value = [double]$dVal
        }
    }

    if ($eVal -ne "BLANK") {
        if ($eVal -eq "NA") {
            $ws.Cells.Item($r, 5).Value = "NA"
        } else {
            $ws.Cells.Item($r, 5).Value = [double]$eVal
        }
    }
}
